$d = $word.ActiveDocument

# --- Locate the paragraph range to replace -------------------------------
# Start paragraph: the list item ending "...par exemple" (the one whose
# trailing run we are rewriting).
# End paragraph:   the empty "Paragraphedeliste"-styled paragraph that
# follows "Ensuite, enregistrer le fichier réparé" (both of these
# paragraphs disappear in the rewrite).

$count = $d.Paragraphs.Count

$startIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*essayer de réparer ces saletés de face*par exemple*") {
        $startIdx = $i
        break
    }
}
if ($startIdx -eq -1) {
    throw "Could not locate the 'par exemple' paragraph"
}

$endIdx = -1
for ($i = $startIdx; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Ensuite, enregistrer le fichier réparé*") {
        $endIdx = $i + 1
        break
    }
}
if ($endIdx -eq -1) {
    throw "Could not locate the paragraph after 'Ensuite, enregistrer le fichier réparé'"
}

$startPara = $d.Paragraphs.Item($startIdx)
$endPara = $d.Paragraphs.Item($endIdx)
$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)

# --- Replacement OOXML -----------------------------------------------------
# Paragraph 1 (same list style as before): trims the old run down to
# "...de réparations", adds a new sentence, a manual line break, and the
# explanation about the "enlever les arêtes surnuméraires" tool, ending on
# the (deliberately ungrammatical, proofed) "les deux face".
# Paragraph 2 (plain): just carries the relocated _GoBack bookmark.
# Paragraph 3 (same list style, no numbering): empty - replaces both the
# old "Ensuite, enregistrer..." paragraph and the blank list paragraph
# that followed it.

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/test.xml"><pkg:xmlData>' +
  '<w:p ' + $ns + '>' +
    '<w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Malheureusement, ça ne marche pas souvent comme prévu, et il faut à tâtons, à l’aide des autres icônes de cet onglet, essayer de réparer ces saletés de face</w:t></w:r>' +
    '<w:r><w:t>s, en testant les outils possibles de réparations</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">. </w:t></w:r>' +
    '<w:r><w:br/><w:t>Souvent, l’outil « </w:t></w:r>' +
    '<w:r><w:t>enlever les arêtes surnuméraires</w:t></w:r>' +
    '<w:r><w:t> »</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">fait combiner les deux faces car ils enlèvent les arêtes entre </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>les deux face</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
  '</w:p>' +
  '<w:p ' + $ns + '>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
  '</w:p>' +
  '<w:p ' + $ns + '>' +
    '<w:pPr><w:pStyle w:val="Paragraphedeliste"/></w:pPr>' +
  '</w:p>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$rng.InsertXML($xml)
